$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename header cells (D1, E1, F1)
$ws.Range("D1").Value = "t1"
$ws.Range("E1").Value = "t2"
$ws.Range("F1").Value = "effect"

# Fill in previously-blank "id" column values that were left empty
# for rows belonging to the same trial/id as the row above them
$ws.Range("C3").Value = 1
$ws.Range("C5").Value = 2
$ws.Range("C6").Value = 2
$ws.Range("C7").Value = 2

# Add new note cell for row 12
$ws.Range("M12").Value = "effect=HR"
